$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 888.5
$ws.Range("I43").Value = 840.5
$ws.Range("J43").Value = 912.5
$ws.Range("K43").Value = 840.5
$ws.Range("L43").Value = 912.5
$ws.Range("M43").Value = -771.5
$ws.Range("N43").Value = -1050.5
$ws.Range("H64").Value = 3904.125
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 3992.2307
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 3992.2307
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -4488.2307
$ws.Range("H67").Value = 3904.125
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 3992.2307
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 3992.2307
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -5708.2307
$ws.Range("H113").Value = 1743.75
$ws.Range("I113").Value = 1517.6364
$ws.Range("J113").Value = 1862.1904
$ws.Range("K113").Value = 1517.6364
$ws.Range("L113").Value = 1862.1904
$ws.Range("M113").Value = 1736.3636
$ws.Range("N113").Value = -8370.190399999999
$ws.Range("H121").Value = 3994.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3994.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 11983.5
$ws.Range("H125").Value = 1175
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1175
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 10575
$ws.Range("N125").Value = -15495
$ws.Range("H131").Value = 24853.326
$ws.Range("I131").Value = 29655.686
$ws.Range("J131").Value = 3843
$ws.Range("K131").Value = 88967.058
$ws.Range("L131").Value = 11529
$ws.Range("M131").Value = -83927.058
$ws.Range("N131").Value = -21609

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8076.643
$ws.Range("I74").Value = 9304.916999999999
$ws.Range("J74").Value = 707
$ws.Range("K74").Value = 9304.916999999999
$ws.Range("L74").Value = 707
$ws.Range("M74").Value = -8430.916999999999
$ws.Range("H77").Value = 8076.643
$ws.Range("I77").Value = 9304.916999999999
$ws.Range("J77").Value = 707
$ws.Range("K77").Value = 46524.585
$ws.Range("L77").Value = 3535
$ws.Range("M77").Value = -42156.585
$ws.Range("H88").Value = 2116.5386
$ws.Range("I88").Value = 2215.1428
$ws.Range("J88").Value = 2001.5
$ws.Range("K88").Value = 2215.1428
$ws.Range("L88").Value = 2001.5
$ws.Range("M88").Value = -1809.1428
$ws.Range("H91").Value = 2116.5386
$ws.Range("I91").Value = 2215.1428
$ws.Range("J91").Value = 2001.5
$ws.Range("K91").Value = 2215.1428
$ws.Range("L91").Value = 2001.5
$ws.Range("M91").Value = -811.1428000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1738.7142
$ws.Range("I86").Value = 1713.8
$ws.Range("J86").Value = 1801
$ws.Range("K86").Value = 1713.8
$ws.Range("L86").Value = 1801
$ws.Range("M86").Value = -590.8
$ws.Range("N86").Value = -4047
$ws.Range("H89").Value = 1738.7142
$ws.Range("I89").Value = 1713.8
$ws.Range("J89").Value = 1801
$ws.Range("K89").Value = 8569
$ws.Range("L89").Value = 9005
$ws.Range("M89").Value = -2953
$ws.Range("N89").Value = -20237
$ws.Range("H132").Value = 41950
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 41950
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 41950
$ws.Range("N132").Value = -52070
$ws.Range("H134").Value = 2433.9285
$ws.Range("I134").Value = 1838.4166
$ws.Range("J134").Value = 6007
$ws.Range("K134").Value = 5515.2498
$ws.Range("L134").Value = 18021
$ws.Range("M134").Value = -2980.2498
$ws.Range("N134").Value = -23091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3713.5
$ws.Range("I110").Value = 2170.25
$ws.Range("J110").Value = 6800
$ws.Range("K110").Value = 6510.75
$ws.Range("L110").Value = 20400
$ws.Range("M110").Value = -2420.75
$ws.Range("N110").Value = -28580
$ws.Range("H131").Value = 2757.698
$ws.Range("I131").Value = 642.5
$ws.Range("J131").Value = 3133.7334
$ws.Range("K131").Value = 1927.5
$ws.Range("L131").Value = 9401.200199999999
$ws.Range("M131").Value = 3112.5
$ws.Range("N131").Value = -19481.2002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2435.4868
$ws.Range("I132").Value = 2200.1304
$ws.Range("J132").Value = 4755.4287
$ws.Range("K132").Value = 6600.3912
$ws.Range("L132").Value = 14266.2861
$ws.Range("M132").Value = -4070.3912
$ws.Range("N132").Value = -19326.2861
$ws.Range("H136").Value = 2580.9375
$ws.Range("I136").Value = 1865.8334
$ws.Range("J136").Value = 4726.25
$ws.Range("K136").Value = 5597.5002
$ws.Range("L136").Value = 14178.75
$ws.Range("M136").Value = -3047.5002
$ws.Range("N136").Value = -19278.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18313
$ws.Range("I45").Value = 18000
$ws.Range("J45").Value = 18626
$ws.Range("K45").Value = 18000
$ws.Range("L45").Value = 18626
$ws.Range("M45").Value = -17509
$ws.Range("N45").Value = -19608
$ws.Range("H62").Value = 7151885.5
$ws.Range("I62").Value = 14301429
$ws.Range("J62").Value = 2342.7144
$ws.Range("K62").Value = 14301429
$ws.Range("L62").Value = 2342.7144
$ws.Range("M62").Value = -14300805
$ws.Range("N62").Value = -3590.7144
$ws.Range("H65").Value = 7151885.5
$ws.Range("I65").Value = 14301429
$ws.Range("J65").Value = 2342.7144
$ws.Range("K65").Value = 71507145
$ws.Range("L65").Value = 11713.572
$ws.Range("M65").Value = -71504025
$ws.Range("N65").Value = -17953.572
$ws.Range("H100").Value = 645.2727
$ws.Range("I100").Value = 337.35715
$ws.Range("J100").Value = 1184.125
$ws.Range("K100").Value = 674.7143
$ws.Range("L100").Value = 2368.25
$ws.Range("M100").Value = -133.7143
$ws.Range("N100").Value = -3450.25
$ws.Range("H126").Value = 2437.4666
$ws.Range("I126").Value = 2437.4666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7312.399800000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4842.399800000001
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 10034.615
$ws.Range("I136").Value = 12216.111
$ws.Range("J136").Value = 5126.25
$ws.Range("K136").Value = 36648.333
$ws.Range("L136").Value = 15378.75
$ws.Range("M136").Value = -34098.333
$ws.Range("N136").Value = -20478.75

